$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = "Teste Jogo PC"
$ws.Range("B82").Value = "Completo"
$ws.Range("C82").Value = "PC"
$ws.Range("D82").Value = "Concluído"
